# Weekly crime-stat refresh: cs-en-us-076pct.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title block: new volume/issue number and report week ---
$ws.Range("A8").Value = "Volume 32   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/1/2025  Through  12/7/2025"

# --- Step 1: convert cells whose data type flips (number <-> N/A text) ---
# Donor cells are in row 14, which this week's update never touches, so their
# style+type (text 0 / text ***.* / formatted number) is a safe template to clone.
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("I14").Copy($ws.Range("D18"))
$ws.Range("L14").Copy($ws.Range("E18"))
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("I14").Copy($ws.Range("C22"))
$ws.Range("I14").Copy($ws.Range("F22"))
$ws.Range("C14").Copy($ws.Range("C25"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("G28"))
$ws.Range("E14").Copy($ws.Range("H28"))

# --- Step 2: write this week's figures into every numeric cell that changed ---
$ws.Range("F16").Value = 4
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 49
$ws.Range("K16").Value = -14.035087719298
$ws.Range("L16").Value = -23.4375
$ws.Range("M16").Value = -49.484536082474
$ws.Range("N16").Value = -87.037037037037
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 7.692307692307
$ws.Range("I17").Value = 106
$ws.Range("J17").Value = 116
$ws.Range("K17").Value = -8.620689655172
$ws.Range("L17").Value = 10.416666666666
$ws.Range("M17").Value = 23.255813953488
$ws.Range("N17").Value = -63.448275862069
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 133.333333333333
$ws.Range("I18").Value = 100
$ws.Range("J18").Value = 82
$ws.Range("K18").Value = 21.951219512195
$ws.Range("L18").Value = 23.456790123456
$ws.Range("M18").Value = 2.040816326530
$ws.Range("N18").Value = -75.369458128078
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -70
$ws.Range("F19").Value = 14
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = -36.363636363636
$ws.Range("I19").Value = 186
$ws.Range("J19").Value = 202
$ws.Range("K19").Value = -7.920792079207
$ws.Range("L19").Value = 18.471337579617
$ws.Range("M19").Value = -12.676056338028
$ws.Range("N19").Value = -17.333333333333
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 36
$ws.Range("K20").Value = -37.931034482758
$ws.Range("L20").Value = -41.935483870967
$ws.Range("M20").Value = -37.931034482758
$ws.Range("N20").Value = -90.551181102362
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -35.714285714285
$ws.Range("F21").Value = 43
$ws.Range("G21").Value = 44
$ws.Range("H21").Value = -2.272727272727
$ws.Range("I21").Value = 482
$ws.Range("J21").Value = 520
$ws.Range("K21").Value = -7.307692307692
$ws.Range("L21").Value = 3.433476394849
$ws.Range("M21").Value = -13.153153153153
$ws.Range("N21").Value = -71.647058823529
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = -46.153846153846
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 18.181818181818
$ws.Range("I23").Value = 109
$ws.Range("J23").Value = 104
$ws.Range("K23").Value = 4.807692307692
$ws.Range("L23").Value = 6.862745098039
$ws.Range("M23").Value = 45.333333333333
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = -37.5
$ws.Range("F24").Value = 20
$ws.Range("G24").Value = 38
$ws.Range("H24").Value = -47.368421052631
$ws.Range("I24").Value = 424
$ws.Range("J24").Value = 505
$ws.Range("K24").Value = -16.039603960396
$ws.Range("L24").Value = -24.285714285714
$ws.Range("M24").Value = 1.923076923076
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 3
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -82.352941176470
$ws.Range("J25").Value = 245
$ws.Range("K25").Value = -36.734693877551
$ws.Range("L25").Value = -46.180555555555
$ws.Range("C26").Value = 6
$ws.Range("F26").Value = 16
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 23.076923076923
$ws.Range("I26").Value = 182
$ws.Range("K26").Value = 16.666666666666
$ws.Range("L26").Value = 25.517241379310
$ws.Range("M26").Value = -33.088235294117
$ws.Range("L29").Value = -55.555555555555
$ws.Range("L30").Value = -66.666666666666
